# Update "想去人数" (interested-attendee count, column F) figures across the
# four sheets of the 广州-漫展信息 workbook, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 147
$ws.Range("F5").Value = 1300
$ws.Range("F6").Value = 222
$ws.Range("F7").Value = 2480
$ws.Range("F8").Value = 883
$ws.Range("F9").Value = 18576
$ws.Range("F10").Value = 51
$ws.Range("F11").Value = 1889
$ws.Range("F12").Value = 659
$ws.Range("F13").Value = 599
$ws.Range("F14").Value = 325
$ws.Range("F15").Value = 599
$ws.Range("F16").Value = 197
$ws.Range("F17").Value = 199

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 165
$ws.Range("F9").Value = 123
$ws.Range("F15").Value = 68
$ws.Range("F19").Value = 3
$ws.Range("F23").Value = 32

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5875
$ws.Range("F3").Value = 563
$ws.Range("F4").Value = 552

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5875
$ws.Range("F4").Value = 563
$ws.Range("F5").Value = 552
$ws.Range("F8").Value = 147
$ws.Range("F10").Value = 1300
$ws.Range("F12").Value = 222
$ws.Range("F13").Value = 165
$ws.Range("F15").Value = 2480
$ws.Range("F16").Value = 883
$ws.Range("F17").Value = 18576
$ws.Range("F19").Value = 51
$ws.Range("F21").Value = 123
$ws.Range("F22").Value = 123
$ws.Range("F23").Value = 1889
$ws.Range("F24").Value = 659
$ws.Range("F26").Value = 599
$ws.Range("F27").Value = 325
$ws.Range("F28").Value = 599
$ws.Range("F29").Value = 197
$ws.Range("F30").Value = 199
$ws.Range("F37").Value = 68
$ws.Range("F48").Value = 32
